# Re-sort the Maltaspor fantasy-roster sheet: three new players (Cade
# Cunningham, Carlton Carrington, Derrick White) are inserted near the top
# of the roster and the remaining players are reordered, while each
# player keeps their own correct position ("Pozisyon"/B) and team
# ("Takim"/C). Column A (player name), B (position) and C (team) are all
# rewritten together per row so every row stays an internally-consistent
# player/position/team triple, matching the final layout in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2;  A="Cade Cunningham";     B="PG,SG";    C="Detroit Pistons"}
    @{Row=3;  A="Carlton Carrington";  B="PG,SG";    C="Washington Wizards"}
    @{Row=4;  A="Derrick White";       B="PG,SG";    C="Boston Celtics"}
    @{Row=5;  A="Malik Monk";          B="PG,SG,SF"; C="Sacramento Kings"}
    @{Row=6;  A="Ausar Thompson";      B="SF,PF";    C="Detroit Pistons"}
    @{Row=7;  A="Naz Reid";            B="PF,C";     C="Minnesota Timberwolves"}
    @{Row=8;  A="Keon Ellis";          B="SG,SF";    C="Sacramento Kings"}
    @{Row=9;  A="Isaiah Hartenstein";  B="C";        C="Oklahoma City Thunder"}
    @{Row=10; A="Precious Achiuwa";    B="PF,C";     C="New York Knicks"}
    @{Row=11; A="Santi Aldama";        B="PF,C";     C="Memphis Grizzlies"}
    @{Row=12; A="Coby White";          B="PG,SG";    C="Chicago Bulls"}
    @{Row=13; A="Onyeka Okongwu";      B="PF,C";     C="Atlanta Hawks"}
    @{Row=14; A="Devin Vassell";       B="SG,SF";    C="San Antonio Spurs"}
    @{Row=15; A="Anthony Davis";       B="PF,C";     C="Dallas Mavericks"}
    @{Row=16; A="LaMelo Ball";         B="PG,SG";    C="Charlotte Hornets"}
    @{Row=17; A="Damian Lillard";      B="PG";       C="Milwaukee Bucks"}
    @{Row=18; A="Collin Sexton";       B="PG,SG";    C="Utah Jazz"}
    @{Row=19; A="Andrew Wiggins";      B="SF,PF";    C="Miami Heat"}
)

foreach ($u in $updates) {
    $ws.Range("A" + $u.Row).Value = $u.A
    $ws.Range("B" + $u.Row).Value = $u.B
    $ws.Range("C" + $u.Row).Value = $u.C
}
